# Generate Report for Handoff
# Update the f9b94de5 report row (now "Ready for handoff") across the
# Overview, zh-cn and de-de sheets, plus refresh the relevant handoff
# timestamps for that entry.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D2").Value = "2016-07-15 04:07:01"
$ws.Range("D3").Value = "2016-07-15 04:07:01"

# --- zh-cn sheet ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E2").Value = "2016-03-15 04:06:54"
$ws.Range("E3").Value = "2016-03-15 04:06:54"

# --- de-de sheet ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E2").Value = "2016-03-15 04:07:01"
$ws.Range("E3").Value = "2016-03-15 04:07:01"
